$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 53,20
$data[0,0] = 'index'
$data[0,1] = 'lgbm_mae'
$data[0,2] = 'consensus_mae'
$data[0,3] = 'lgbm_mse'
$data[0,4] = 'consensus_mse'
$data[0,5] = 'lgbm_r2'
$data[0,6] = 'consensus_r2'
$data[0,7] = 'consensus_r2_org'
$data[0,8] = 'len'
$data[0,9] = 'consensus_mae_eps'
$data[0,10] = 'lgbm_mse_org'
$data[0,11] = 'consensus_mape_eps'
$data[0,12] = 'consensus_mse_org'
$data[0,13] = 'lgbm_mape_eps'
$data[0,14] = 'lgbm_mae_org'
$data[0,15] = 'lgbm_medae_org'
$data[0,16] = 'consensus_mae_org'
$data[0,17] = 'lgbm_mae_eps'
$data[0,18] = 'lgbm_r2_org'
$data[0,19] = 'consensus_medae_org'
$data[1,0] = 'ibes_1|fwdepsqcut-46|dense2｜compare large space'
$data[1,1] = 0.00966153047961
$data[1,2] = 0.008733878054388279
$data[1,3] = 0.0002005912151791027
$data[1,4] = 0.0001961492659581503
$data[1,5] = 0.1083032380727689
$data[1,6] = 0.1280492261183164
$data[1,7] = 0.2582640092197721
$data[1,8] = 14156
$data[1,10] = 0.0004250480609163007
$data[1,12] = 0.0003384092866814527
$data[1,14] = 0.01071193221870467
$data[1,15] = 0.005638666708596741
$data[1,16] = 0.009674189396799987
$data[1,18] = 0.06836645151013176
$data[1,19] = 0.004818030126325816
$data[2,0] = 'ibes_qoq_2|fwdepsqcut|q_2｜ibes_qoq_tune10_ind3'
$data[2,1] = 0.003127879691820042
$data[2,2] = 0.002538761048289119
$data[2,3] = 0.00002076719194527133
$data[2,4] = 0.00002103311528121684
$data[2,5] = 0.01355620914095823
$data[2,6] = 0.0009248228524665336
$data[2,7] = 0.2133795918767654
$data[2,8] = 11611
$data[2,10] = 0.0001978384306784256
$data[2,12] = 0.0001631540263502446
$data[2,14] = 0.00502885184437347
$data[2,15] = 0.001779221775169934
$data[2,16] = 0.004580501696160511
$data[2,18] = 0.04615441883950755
$data[2,19] = 0.001592932311570713
$data[3,0] = 'ibes_2|ni_depthwise|xgb ind4 -sample_type industry -x_type ni_sp500'
$data[3,1] = 0.005400048253033386
$data[3,2] = 0.005699109365498782
$data[3,3] = 0.0000731390637792532
$data[3,4] = 0.00009470932808225897
$data[3,5] = 0.4252252477925591
$data[3,6] = 0.2557119579146814
$data[3,7] = 0.4208397194991282
$data[3,8] = 6771
$data[3,10] = 0.0001137832146297872
$data[3,12] = 0.0001037297229528328
$data[3,14] = 0.005606434133779571
$data[3,15] = 0.003094044369724517
$data[3,16] = 0.005813234561065725
$data[3,18] = 0.36470746642943
$data[3,19] = 0.003193202088483064
$data[4,0] = 'ibes_2|ni_depthwise|xgb ind_all_tuning -sample_type industry -x_type ni_sp500'
$data[4,1] = 0.005308188556537048
$data[4,2] = 0.005639127357765789
$data[4,3] = 0.00007028985097208806
$data[4,4] = 0.00009102681873913914
$data[4,5] = 0.4089911478654851
$data[4,6] = 0.2346312460124846
$data[4,7] = 0.4208397194991282
$data[4,8] = 6771
$data[4,10] = 0.0001152638129503358
$data[4,12] = 0.0001037297229528328
$data[4,14] = 0.005632685300523813
$data[4,15] = 0.003171635367762127
$data[4,16] = 0.005813234561065725
$data[4,18] = 0.3564407544954944
$data[4,19] = 0.003193202088483064
$data[5,0] = 'ibes_1|fwdepsqcut-49|dense2｜hyperopt_compare3 -code 0 -exclude_fwd False'
$data[5,1] = 0.008311421494257554
$data[5,2] = 0.008853442839112767
$data[5,3] = 0.0001712304539584277
$data[5,4] = 0.0002150273861825077
$data[5,5] = 0.3115448762458137
$data[5,6] = 0.1354534059650468
$data[5,7] = 0.2582640092197723
$data[5,8] = 42468
$data[5,10] = 0.000333070102780211
$data[5,12] = 0.0003384092866814527
$data[5,14] = 0.009106885821888167
$data[5,15] = 0.004616722796530505
$data[5,16] = 0.009674189396799987
$data[5,18] = 0.269966598412229
$data[5,19] = 0.004818030126325817
$data[6,0] = 'ibes_2|ni_depthwise|xgb ind_all_tuning -sample_type industry -x_type ni'
$data[6,1] = 0.007944638673415367
$data[6,2] = 0.008722720041403264
$data[6,3] = 0.0001442016437287739
$data[6,4] = 0.0001985649052123504
$data[6,5] = 0.3807705094595015
$data[6,6] = 0.147324247391146
$data[6,7] = 0.2582640092197721
$data[6,8] = 14156
$data[6,10] = 0.0003285627388197843
$data[6,12] = 0.0003384092866814527
$data[6,14] = 0.008950010380536522
$data[6,15] = 0.004615286290830784
$data[6,16] = 0.009674189396799985
$data[6,18] = 0.2798459788092014
$data[6,19] = 0.004818030126325817
$data[7,0] = 'ibes_2|fwdepsqcut|ibes_new industry_only ws -indi space3'
$data[7,1] = 0.008392539580425228
$data[7,2] = 0.008722720041403264
$data[7,3] = 0.0001662926371782633
$data[7,4] = 0.0001985649052123504
$data[7,5] = 0.2859075504422648
$data[7,6] = 0.147324247391146
$data[7,7] = 0.2582640092197721
$data[7,8] = 14156
$data[7,10] = 0.0003713407546776094
$data[7,12] = 0.0003384092866814527
$data[7,14] = 0.00941140169945745
$data[7,15] = 0.004590555422836876
$data[7,16] = 0.009674189396799985
$data[7,18] = 0.1860837943045432
$data[7,19] = 0.004818030126325816
$data[8,0] = 'ibes_1|fwdepsqcut-46|dense2｜sp_fix_space -best_col 0 -code 0 -exclude_fwd True'
$data[8,1] = 0.007970745265560365
$data[8,2] = 0.00654101945735325
$data[8,3] = 0.0001659654171746955
$data[8,4] = 0.0001252256324013319
$data[8,5] = -0.0277890197779842
$data[8,6] = 0.224503922757791
$data[8,7] = 0.3127047460379806
$data[8,8] = 24768
$data[8,10] = 0.0002632414236722758
$data[8,12] = 0.0001796235360933539
$data[8,14] = 0.008363098231375462
$data[8,15] = 0.004275415213561056
$data[8,16] = 0.006875190708555331
$data[8,18] = -0.00724317687482956
$data[8,19] = 0.003559726400188349
$data[9,0] = 'ibes_2|fwdepsqcut|cnn_rnn｜industry_exclude'
$data[9,1] = 0.009885890648801
$data[9,2] = 0.009614451672116909
$data[9,3] = 0.0002810826792241289
$data[9,4] = 0.0002968767906711129
$data[9,5] = 0.2059318010033613
$data[9,6] = 0.1613128950427474
$data[9,7] = 0.2582640092197721
$data[9,8] = 14156
$data[9,10] = 0.0003847976244374441
$data[9,12] = 0.0003384092866814527
$data[9,14] = 0.009901025034115019
$data[9,15] = 0.005167135423672592
$data[9,16] = 0.009674189396799985
$data[9,18] = 0.1565886089862186
$data[9,19] = 0.004818030126325816
$data[10,0] = 'ibes_qoq_2|fwdepsqcut|q_2｜ibes_qoq_tune10_ind2'
$data[10,1] = 0.003125769257313045
$data[10,2] = 0.00254141838673398
$data[10,3] = 0.00002166040204619485
$data[10,4] = 0.00002106152912726658
$data[10,5] = -0.02738763561426372
$data[10,6] = 0.001017868166272584
$data[10,7] = 0.212555477012045
$data[10,8] = 11542
$data[10,10] = 0.000209188838038313
$data[10,12] = 0.000163888349337522
$data[10,14] = 0.00512210974064752
$data[10,15] = 0.001602037466360856
$data[10,16] = 0.004589584286176587
$data[10,18] = -0.005102592401121431
$data[10,19] = 0.001597174444307401
$data[11,0] = 'ibes_2|ni|ibes_new industry_all x -indi space_sp500'
$data[11,1] = 0.005058551760597552
$data[11,2] = 0.005699109365498782
$data[11,3] = 0.00007085921366241217
$data[11,4] = 0.00009470932808225897
$data[11,5] = 0.4431418058979847
$data[11,6] = 0.2557119579146814
$data[11,7] = 0.4208397194991282
$data[11,8] = 6771
$data[11,10] = 0.000112952017849086
$data[11,12] = 0.0001037297229528328
$data[11,14] = 0.00526064949972486
$data[11,15] = 0.002682213716882752
$data[11,16] = 0.005813234561065725
$data[11,18] = 0.3693483364419842
$data[11,19] = 0.003193202088483064
$data[12,0] = 'ibes_1|fwdepsqcut|sp500_entire_sp500'
$data[12,1] = 0.005806826536720938
$data[12,2] = 0.005670204101840715
$data[12,3] = 0.00008986367404537444
$data[12,4] = 0.00009344840654260834
$data[12,5] = 0.2696481973337195
$data[12,6] = 0.2405138906269937
$data[12,7] = 0.4208397194991282
$data[12,8] = 6771
$data[12,10] = 0.0001410934357507483
$data[12,12] = 0.0001037297229528328
$data[12,14] = 0.006002887168865449
$data[12,15] = 0.003030311329129174
$data[12,16] = 0.005813234561065723
$data[12,18] = 0.2122246979933391
$data[12,19] = 0.003193202088483064
$data[13,0] = 'ibes_qoq_2|fwdepsqcut|q_2｜ibes_qoq'
$data[13,1] = 0.003208711377129972
$data[13,2] = 0.002538761048289119
$data[13,3] = 0.00002210896317387241
$data[13,4] = 0.00002103311528121684
$data[13,5] = -0.05017806464506647
$data[13,6] = 0.0009248228524665336
$data[13,7] = 0.2133795918767654
$data[13,8] = 11611
$data[13,10] = 0.0002098592232396584
$data[13,12] = 0.0001631540263502446
$data[13,14] = 0.005213165544267774
$data[13,15] = 0.00173711317318318
$data[13,16] = 0.004580501696160511
$data[13,18] = -0.01180186309853548
$data[13,19] = 0.001592932311570713
$data[14,0] = 'ibes_qoq_2|ni|q_2｜ibes_qoq'
$data[14,1] = 0.003491208235511574
$data[14,2] = 0.001329144946829925
$data[14,3] = 0.0000185690054965323
$data[14,4] = 0.00000916857531699537
$data[14,5] = -0.009024320017351872
$data[14,6] = 0.501787240222058
$data[14,7] = 0.6230156929756394
$data[14,8] = 49
$data[14,10] = 0.0003123469992588732
$data[14,12] = 0.0001158352549329261
$data[14,14] = 0.007217175147339313
$data[14,15] = 0.002972419699024923
$data[14,16] = 0.004529835290885987
$data[14,18] = -0.01652918306199047
$data[14,19] = 0.00129198945232948
$data[15,0] = 'ibes_1|fwdepsqcut|cnn_rnn｜without ibes -2'
$data[15,1] = 0.01042798137117382
$data[15,2] = 0.009812042985245533
$data[15,3] = 0.0002839828498694114
$data[15,4] = 0.000297330738550053
$data[15,5] = 0.1884050333044567
$data[15,6] = 0.1502580843805948
$data[15,7] = 0.2597560818810524
$data[15,8] = 12741
$data[15,10] = 0.0004007814960332574
$data[15,12] = 0.0003441387227458278
$data[15,14] = 0.01033009227796955
$data[15,15] = 0.005554310648909318
$data[15,16] = 0.00971495424881711
$data[15,18] = 0.1379172254546043
$data[15,19] = 0.004797770485226469
$data[16,0] = 'ibes_1|ni-industry_code|cnn_rnn｜without ibes -2'
$data[16,1] = 0.00997695293896445
$data[16,2] = 0.009631259951277263
$data[16,3] = 0.0002746516022984678
$data[16,4] = 0.000318358280400409
$data[16,5] = 0.272470124651383
$data[16,6] = 0.156694670201814
$data[16,7] = 0.4558883625799528
$data[16,8] = 1284
$data[16,10] = 0.0003741973334738424
$data[16,12] = 0.0002640644396388745
$data[16,14] = 0.009685349739034687
$data[16,15] = 0.005106524251936465
$data[16,16] = 0.009204693368560462
$data[16,18] = 0.2289566739349258
$data[16,19] = 0.004252957223118693
$data[17,0] = 'ibes_1|ni-sector_code|cnn_rnn｜without ibes -2'
$data[17,1] = 0.0102733372090916
$data[17,2] = 0.00987889837479569
$data[17,3] = 0.0002816309717823605
$data[17,4] = 0.0003034120436381377
$data[17,5] = 0.204771852061556
$data[17,6] = 0.1432696624324651
$data[17,7] = 0.2595879210340796
$data[17,8] = 12051
$data[17,10] = 0.0004008200453677453
$data[17,12] = 0.0003482779314587717
$data[17,14] = 0.01017444484578581
$data[17,15] = 0.005371111699895803
$data[17,16] = 0.009769435123820693
$data[17,18] = 0.1478874304814285
$data[17,19] = 0.004769747765218425
$data[18,0] = 'ibes_1|fwdepsqcut-46|dense2｜all x 0 -fix space'
$data[18,1] = 0.009164293356549506
$data[18,2] = 0.008733878054388279
$data[18,3] = 0.0001865208677798922
$data[18,4] = 0.0001961492659581503
$data[18,5] = 0.1708507589293765
$data[18,6] = 0.1280492261183164
$data[18,7] = 0.2582640092197721
$data[18,8] = 14156
$data[18,10] = 0.0004006378783650764
$data[18,12] = 0.0003384092866814527
$data[18,14] = 0.0102119158074488
$data[18,15] = 0.005239736253559311
$data[18,16] = 0.009674189396799987
$data[18,18] = 0.1218694481841025
$data[18,19] = 0.004818030126325817
$data[19,0] = 'ibes_2|ni|ibes_new industry_all x -mse_sp500'
$data[19,1] = 0.005791630591612971
$data[19,2] = 0.005699109365498782
$data[19,3] = 0.00008698446620952695
$data[19,4] = 0.00009470932808225897
$data[19,5] = 0.3164189910555092
$data[19,6] = 0.2557119579146814
$data[19,7] = 0.4208397194991282
$data[19,8] = 6771
$data[19,10] = 0.0001324165619288979
$data[19,12] = 0.0001037297229528328
$data[19,14] = 0.006002882908757912
$data[19,15] = 0.003330775568204595
$data[19,16] = 0.005813234561065725
$data[19,18] = 0.2606707993949464
$data[19,19] = 0.003193202088483064
$data[20,0] = 'ibes_2|fwdepsqcut|ibes_industry -sp500_sp500'
$data[20,1] = 0.005611881407010991
$data[20,2] = 0.005639127357765789
$data[20,3] = 0.00008431607544410066
$data[20,4] = 0.00009102681873913914
$data[20,5] = 0.291056300795216
$data[20,6] = 0.2346312460124846
$data[20,7] = 0.4208397194991282
$data[20,8] = 6771
$data[20,10] = 0.0001424514975481683
$data[20,12] = 0.0001037297229528328
$data[20,14] = 0.00595976388177965
$data[20,15] = 0.002972129451884979
$data[20,16] = 0.005813234561065725
$data[20,18] = 0.2046421514566161
$data[20,19] = 0.003193202088483064
$data[21,0] = 'ibes_2|fwdepsqcut|ibes_new industry_only ws -indi space3_sp500'
$data[21,1] = 0.005713567326092179
$data[21,2] = 0.005699109365498782
$data[21,3] = 0.00008903433264960941
$data[21,4] = 0.00009470932808225897
$data[21,5] = 0.3003097955821736
$data[21,6] = 0.2557119579146814
$data[21,7] = 0.4208397194991282
$data[21,8] = 6771
$data[21,10] = 0.00013855240761432
$data[21,12] = 0.0001037297229528328
$data[21,14] = 0.00592641693601841
$data[21,15] = 0.003005755874836842
$data[21,16] = 0.005813234561065725
$data[21,18] = 0.2264121702660996
$data[21,19] = 0.003193202088483064
$data[22,0] = 'ibes_2|ni_depthwise|xgb ind3 -sample_type industry -x_type ni'
$data[22,1] = 0.007927365294981236
$data[22,2] = 0.008722720041403264
$data[22,3] = 0.000143672851475712
$data[22,4] = 0.0001985649052123504
$data[22,5] = 0.3830412447229715
$data[22,6] = 0.147324247391146
$data[22,7] = 0.2582640092197721
$data[22,8] = 14156
$data[22,10] = 0.0003279955215444224
$data[22,12] = 0.0003384092866814527
$data[22,14] = 0.008942022948316988
$data[22,15] = 0.004617426367895965
$data[22,16] = 0.009674189396799985
$data[22,18] = 0.2810892232598902
$data[22,19] = 0.004818030126325817
$data[23,0] = 'ibes_qoq_1|fwdepsqcut|q_1｜ibes_qoq_tune10'
$data[23,1] = 0.003031636244261009
$data[23,2] = 0.002499571894821071
$data[23,3] = 0.0000191501667802569
$data[23,4] = 0.00001967007435869721
$data[23,5] = -0.01488433510101705
$data[23,6] = -0.04243741404356438
$data[23,7] = 0.2133795918767654
$data[23,8] = 11611
$data[23,10] = 0.0002083353057278959
$data[23,12] = 0.0001631540263502446
$data[23,14] = 0.005089153979501482
$data[23,15] = 0.001554417436991534
$data[23,16] = 0.004580501696160511
$data[23,18] = -0.004454544482717981
$data[23,19] = 0.001592932311570713
$data[24,0] = 'ibes_2|ni_depthwise|xgb ind4 -sample_type industry -x_type ni'
$data[24,1] = 0.00784507820846319
$data[24,2] = 0.008722720041403264
$data[24,3] = 0.0001417160098566256
$data[24,4] = 0.0001985649052123504
$data[24,5] = 0.3914442976114285
$data[24,6] = 0.147324247391146
$data[24,7] = 0.2582640092197721
$data[24,8] = 14156
$data[24,10] = 0.0003254854012574867
$data[24,12] = 0.0003384092866814527
$data[24,14] = 0.008855143980453176
$data[24,15] = 0.00456770441607247
$data[24,16] = 0.009674189396799985
$data[24,18] = 0.2865909829079945
$data[24,19] = 0.004818030126325816
$data[25,0] = 'ibes_qoq_2|fwdepsqcut|q_2｜ibes_qoq_tune10_ind'
$data[25,1] = 0.003103990424644728
$data[25,2] = 0.00254141838673398
$data[25,3] = 0.00002156002471405648
$data[25,4] = 0.00002106152912726658
$data[25,5] = -0.02262657763782516
$data[25,6] = 0.001017868166272584
$data[25,7] = 0.212555477012045
$data[25,8] = 11542
$data[25,10] = 0.0002091701444985722
$data[25,12] = 0.000163888349337522
$data[25,14] = 0.00510134746390362
$data[25,15] = 0.001583264354559398
$data[25,16] = 0.004589584286176587
$data[25,18] = -0.005012774390606101
$data[25,19] = 0.001597174444307401
$data[26,0] = 'ibes_1|fwdepsqcut|ibes_entire_only ws -smaller space_sp500'
$data[26,1] = 0.005895505445859815
$data[26,2] = 0.005670204101840715
$data[26,3] = 0.00009166134936686534
$data[26,4] = 0.00009344840654260834
$data[26,5] = 0.255037895389058
$data[26,6] = 0.2405138906269937
$data[26,7] = 0.4208397194991282
$data[26,8] = 6771
$data[26,10] = 0.0001427653829362463
$data[26,12] = 0.0001037297229528328
$data[26,14] = 0.006096382291329448
$data[26,15] = 0.003125364361436701
$data[26,16] = 0.005813234561065723
$data[26,18] = 0.2028896166553128
$data[26,19] = 0.003193202088483064
$data[27,0] = 'ibes_1|fwdepsqcut-industry_code|ibes_entire_only ws -smaller space_sp500'
$data[27,1] = 0.005903838435270049
$data[27,2] = 0.005670204101840715
$data[27,3] = 0.00009181732856753832
$data[27,4] = 0.00009344840654260834
$data[27,5] = 0.2537702008328302
$data[27,6] = 0.2405138906269937
$data[27,7] = 0.4208397194991282
$data[27,8] = 6771
$data[27,10] = 0.0001433490956550128
$data[27,12] = 0.0001037297229528328
$data[27,14] = 0.00610593491672969
$data[27,15] = 0.003163476470023769
$data[27,16] = 0.005813234561065723
$data[27,18] = 0.1996305390032257
$data[27,19] = 0.003193202088483064
$data[28,0] = 'ibes_1|fwdepsqcut-sector_code|ibes_entire_only ws -smaller space_sp500'
$data[28,1] = 0.005871671098387516
$data[28,2] = 0.005670204101840715
$data[28,3] = 0.00009132537185586264
$data[28,4] = 0.00009344840654260834
$data[28,5] = 0.2577684957503591
$data[28,6] = 0.2405138906269937
$data[28,7] = 0.4208397194991282
$data[28,8] = 6771
$data[28,10] = 0.0001429245641932233
$data[28,12] = 0.0001037297229528328
$data[28,14] = 0.006067498374332903
$data[28,15] = 0.003090643964573058
$data[28,16] = 0.005813234561065723
$data[28,18] = 0.2020008505541708
$data[28,19] = 0.003193202088483064
$data[29,0] = 'ibes_qoq_1|fwdepsqcut|q_1｜ibes_qoq'
$data[29,1] = 0.00304699828992123
$data[29,2] = 0.002499571894821071
$data[29,3] = 0.00001903121423403654
$data[29,4] = 0.00001967007435869721
$data[29,5] = -0.008580312939503454
$data[29,6] = -0.04243741404356438
$data[29,7] = 0.2133795918767654
$data[29,8] = 11611
$data[29,10] = 0.000207986162951978
$data[29,12] = 0.0001631540263502446
$data[29,14] = 0.005107549406927882
$data[29,15] = 0.001554596452837835
$data[29,16] = 0.004580501696160511
$data[29,18] = -0.002771209789547591
$data[29,19] = 0.001592932311570713
$data[30,0] = 'ibes_1|fwdepsqcut-46|dense2｜new with indi code -fix space_sp500'
$data[30,1] = 0.006549373422714545
$data[30,2] = 0.005670204101840715
$data[30,3] = 0.0001064697656049724
$data[30,4] = 0.00009344840654260834
$data[30,5] = 0.1346849985258249
$data[30,6] = 0.2405138906269937
$data[30,7] = 0.4208397194991282
$data[30,8] = 6771
$data[30,10] = 0.000160150142460433
$data[30,12] = 0.0001037297229528328
$data[30,14] = 0.006776662396495037
$data[30,15] = 0.003715325925925931
$data[30,16] = 0.005813234561065723
$data[30,18] = 0.1058242633905934
$data[30,19] = 0.003193202088483064
$data[31,0] = 'ibes_qoq_1|fwdepsqcut|q_1｜ibes_qoq_filter'
$data[31,1] = 0.002783525716284511
$data[31,2] = 0.002457544020376962
$data[31,3] = 0.00001713865210259099
$data[31,4] = 0.0000189022549270271
$data[31,5] = -0.05370835704473098
$data[31,6] = -0.1621371309933592
$data[31,7] = 0.3649062991209757
$data[31,8] = 4010
$data[31,10] = 0.0002742305932959652
$data[31,12] = 0.0001730864695810989
$data[31,14] = 0.004501420066184821
$data[31,15] = 0.001362431751651251
$data[31,16] = 0.004349183403839366
$data[31,18] = -0.006214540120261569
$data[31,19] = 0.001438312583710056
$data[32,0] = 'ibes_1|fwdepsqcut-46|dense2｜mini_tune15_re -code 0 -exclude_fwd True'
$data[32,1] = 0.00929322862430519
$data[32,2] = 0.008733878054388279
$data[32,3] = 0.000192984991667648
$data[32,4] = 0.0001961492659581503
$data[32,5] = 0.1421155108066621
$data[32,6] = 0.1280492261183164
$data[32,7] = 0.2582640092197721
$data[32,8] = 14156
$data[32,10] = 0.0004001063569688302
$data[32,12] = 0.0003384092866814527
$data[32,14] = 0.01035431191953904
$data[32,15] = 0.005409147625572685
$data[32,16] = 0.009674189396799987
$data[32,18] = 0.1230344532976788
$data[32,19] = 0.004818030126325817
$data[33,0] = 'ibes_1|fwdepsqcut|cnn_rnn｜small_training_True_0'
$data[33,1] = 0.01016347810751514
$data[33,2] = 0.009447153043753395
$data[33,3] = 0.0002829041504470825
$data[33,4] = 0.0002806495860517923
$data[33,5] = 0.1619293810393613
$data[33,6] = 0.1686082656554458
$data[33,7] = 0.3904390743053188
$data[33,8] = 4699
$data[33,10] = 0.0004705186388866096
$data[33,12] = 0.0003220429185286173
$data[33,14] = 0.01022526679210232
$data[33,15] = 0.005350192890365459
$data[33,16] = 0.009384601794708014
$data[33,18] = 0.1094051116331665
$data[33,19] = 0.004731165595531388
$data[34,0] = 'ibes_qoq_1|fwdepsqcut|q_1｜ibes_qoqcut8_entire'
$data[34,1] = 0.003030458698074221
$data[34,2] = 0.002499571894821071
$data[34,3] = 0.0000191281128971015
$data[34,4] = 0.00001967007435869721
$data[34,5] = -0.01371556509502225
$data[34,6] = -0.04243741404356438
$data[34,7] = 0.2133795918767654
$data[34,8] = 11611
$data[34,10] = 0.0002081444394071696
$data[34,12] = 0.0001631540263502446
$data[34,14] = 0.005089488294347637
$data[34,15] = 0.001547174514616301
$data[34,16] = 0.004580501696160511
$data[34,18] = -0.00353431378743374
$data[34,19] = 0.001592932311570713
$data[35,0] = 'ibes_2|ni|ibes_new industry_all x -mse'
$data[35,1] = 0.008384122629515654
$data[35,2] = 0.008722720041403264
$data[35,3] = 0.0001620814063416336
$data[35,4] = 0.0001985649052123504
$data[35,5] = 0.3039913826239528
$data[35,6] = 0.147324247391146
$data[35,7] = 0.2582640092197721
$data[35,8] = 14156
$data[35,10] = 0.0003610907680501377
$data[35,12] = 0.0003384092866814527
$data[35,14] = 0.009406309171688502
$data[35,15] = 0.004884082851211167
$data[35,16] = 0.009674189396799985
$data[35,18] = 0.2085500335179142
$data[35,19] = 0.004818030126325816
$data[36,0] = 'ibes_2|fwdepsqcut|rounding_ind_ex'
$data[36,1] = 0.008519531777817098
$data[36,2] = 0.008722720041403264
$data[36,3] = 0.0001673145877462598
$data[36,4] = 0.0001985649052123504
$data[36,5] = 0.2815190988739301
$data[36,6] = 0.147324247391146
$data[36,7] = 0.2582640092197721
$data[36,8] = 14156
$data[36,10] = 0.0003710863970746162
$data[36,12] = 0.0003384092866814527
$data[36,14] = 0.00953579678993187
$data[36,15] = 0.004784771004405694
$data[36,16] = 0.009674189396799985
$data[36,18] = 0.1866413032030688
$data[36,19] = 0.004818030126325816
$data[37,0] = 'ibes_1|fwdepsqcut|ibes_entire_only ws -smaller space'
$data[37,1] = 0.008532806478705164
$data[37,2] = 0.008733878054388279
$data[37,3] = 0.0001656454412935777
$data[37,4] = 0.0001961492659581503
$data[37,5] = 0.2636491907304709
$data[37,6] = 0.1280492261183164
$data[37,7] = 0.2582640092197721
$data[37,8] = 14156
$data[37,10] = 0.0003744000227958637
$data[37,12] = 0.0003384092866814527
$data[37,14] = 0.009554445336689437
$data[37,15] = 0.004730871750625608
$data[37,16] = 0.009674189396799987
$data[37,18] = 0.179378395374721
$data[37,19] = 0.004818030126325817
$data[38,0] = 'ibes_1|fwdepsqcut-industry_code|ibes_entire_only ws -smaller space'
$data[38,1] = 0.008529233938818522
$data[38,2] = 0.008733878054388279
$data[38,3] = 0.0001653251466070655
$data[38,4] = 0.0001961492659581503
$data[38,5] = 0.2650730104853413
$data[38,6] = 0.1280492261183164
$data[38,7] = 0.2582640092197721
$data[38,8] = 14156
$data[38,10] = 0.0003746361399868205
$data[38,12] = 0.0003384092866814527
$data[38,14] = 0.009555217528480035
$data[38,15] = 0.004764640177213803
$data[38,16] = 0.009674189396799987
$data[38,18] = 0.1788608663781261
$data[38,19] = 0.004818030126325817
$data[39,0] = 'ibes_1|fwdepsqcut-sector_code|ibes_entire_only ws -smaller space'
$data[39,1] = 0.008505250051591662
$data[39,2] = 0.008733878054388279
$data[39,3] = 0.0001646168902522808
$data[39,4] = 0.0001961492659581503
$data[39,5] = 0.268221452941368
$data[39,6] = 0.1280492261183164
$data[39,7] = 0.2582640092197721
$data[39,8] = 14156
$data[39,10] = 0.0003736225585176234
$data[39,12] = 0.0003384092866814527
$data[39,14] = 0.009526287265811475
$data[39,15] = 0.004756476812602748
$data[39,16] = 0.009674189396799987
$data[39,18] = 0.1810824657398455
$data[39,19] = 0.004818030126325817
$data[40,0] = 'ibes_1|fwdepsqcut-46|dense2｜sp_fix_space -best_col 0 -code 0 -exclude_fwd True_sp500'
$data[40,1] = 0.006299187529923155
$data[40,2] = 0.00553167406053466
$data[40,3] = 0.0000900395869891457
$data[40,4] = 0.00008170381435551846
$data[40,5] = 0.1391419988807219
$data[40,6] = 0.2188393498696173
$data[40,7] = 0.4208397194991282
$data[40,8] = 6771
$data[40,10] = 0.0001603751212724098
$data[40,12] = 0.0001037297229528328
$data[40,14] = 0.00674135934653288
$data[40,15] = 0.003667719421830667
$data[40,16] = 0.005813234561065723
$data[40,18] = 0.1045681259196538
$data[40,19] = 0.003193202088483064
$data[41,0] = 'ibes_1|fwdepsqcut_depthwise|xgb tryrun -sample_type entire -x_type fwdepsqcut'
$data[41,1] = 0.00864935396008955
$data[41,2] = 0.008733878054388279
$data[41,3] = 0.0001602711781363747
$data[41,4] = 0.0001961492659581503
$data[41,5] = 0.2875396340419782
$data[41,6] = 0.1280492261183164
$data[41,7] = 0.2582640092197721
$data[41,8] = 14156
$data[41,10] = 0.0003609428129771297
$data[41,12] = 0.0003384092866814527
$data[41,14] = 0.009669808352814117
$data[41,15] = 0.005274143159156646
$data[41,16] = 0.009674189396799987
$data[41,18] = 0.2088743260447083
$data[41,19] = 0.004818030126325816
$data[42,0] = 'ibes_2|fwdepsqcut_depthwise|xgb ind2 -sample_type industry -x_type fwdepsqcut'
$data[42,1] = 0.00869522801124366
$data[42,2] = 0.008722720041403264
$data[42,3] = 0.0001671499749450426
$data[42,4] = 0.0001985649052123504
$data[42,5] = 0.2822259777859754
$data[42,6] = 0.147324247391146
$data[42,7] = 0.2582640092197721
$data[42,8] = 14156
$data[42,10] = 0.0003635971635380399
$data[42,12] = 0.0003384092866814527
$data[42,14] = 0.009719161233546112
$data[42,15] = 0.005173194757001598
$data[42,16] = 0.009674189396799985
$data[42,18] = 0.2030564380000814
$data[42,19] = 0.004818030126325817
$data[43,0] = 'ibes_2|ni|ibes_new industry_all x -indi space'
$data[43,1] = 0.007668256776450208
$data[43,2] = 0.008722720041403264
$data[43,3] = 0.0001436777206558706
$data[43,4] = 0.0001985649052123504
$data[43,5] = 0.3830203355303232
$data[43,6] = 0.147324247391146
$data[43,7] = 0.2582640092197721
$data[43,8] = 14156
$data[43,10] = 0.0003344607981243851
$data[43,12] = 0.0003384092866814527
$data[43,14] = 0.00867861064019221
$data[43,15] = 0.004128330862335852
$data[43,16] = 0.009674189396799985
$data[43,18] = 0.2669184291403396
$data[43,19] = 0.004818030126325817
$data[44,0] = 'ibes_1|ni|cnn_rnn｜small_training_False_0'
$data[44,1] = 0.01020612284563066
$data[44,2] = 0.009760820655321729
$data[44,3] = 0.0002760452274294578
$data[44,4] = 0.0002924361182664438
$data[44,5] = 0.1971347954963554
$data[44,6] = 0.1494626221848404
$data[44,7] = 0.2582640092197721
$data[44,8] = 14156
$data[44,10] = 0.0003886501570187028
$data[44,12] = 0.0003384092866814527
$data[44,14] = 0.01011738460212881
$data[44,15] = 0.005394960099510664
$data[44,16] = 0.009674189396799987
$data[44,18] = 0.1481445083553075
$data[44,19] = 0.004818030126325816
$data[45,0] = 'ibes_2|fwdepsqcut-46|dense2｜new industry model -fix space'
$data[45,1] = 0.00958912348219286
$data[45,2] = 0.008722720041403264
$data[45,3] = 0.00021842684931279
$data[45,4] = 0.0001985649052123504
$data[45,5] = 0.06203325341612431
$data[45,6] = 0.147324247391146
$data[45,7] = 0.2582640092197721
$data[45,8] = 14156
$data[45,10] = 0.0004313797729252501
$data[45,12] = 0.0003384092866814527
$data[45,14] = 0.01064772150300563
$data[45,15] = 0.005560650653515711
$data[45,16] = 0.009674189396799985
$data[45,18] = 0.05448840836791147
$data[45,19] = 0.004818030126325817
$data[46,0] = 'ibes_2|fwdepsqcut|ibes_industry -sp500'
$data[46,1] = 0.005414655366858831
$data[46,2] = 0.005250633999685894
$data[46,3] = 0.00008390458217824665
$data[46,4] = 0.00008336263472823423
$data[46,5] = 0.2548808040617673
$data[46,6] = 0.2596935978056915
$data[46,7] = 0.456874413314975
$data[46,8] = 5153
$data[46,9] = 0.4055123229186881
$data[46,10] = 0.000141190946981845
$data[46,11] = 10.4129840341794
$data[46,12] = 0.00009243123959608208
$data[46,13] = 11.43997477570016
$data[46,14] = 0.005626234126090045
$data[46,15] = 0.002733726006933551
$data[46,16] = 0.005318929375230071
$data[46,17] = 0.4113135413420823
$data[46,18] = 0.1703625716885953
$data[46,19] = 0.002848413552158193
$data[47,0] = 'ibes_1|fwdepsqcut-46|dense2｜large_big_allx -code 0 -exclude_fwd True'
$data[47,1] = 0.00955739351620514
$data[47,2] = 0.008733878054388279
$data[47,3] = 0.0001939926155429641
$data[47,4] = 0.0001961492659581503
$data[47,5] = 0.137636277027368
$data[47,6] = 0.1280492261183164
$data[47,7] = 0.2582640092197721
$data[47,8] = 14156
$data[47,10] = 0.0004139067990685161
$data[47,12] = 0.0003384092866814527
$data[47,14] = 0.01060174085834054
$data[47,15] = 0.005675761589993111
$data[47,16] = 0.009674189396799987
$data[47,18] = 0.09278621544819221
$data[47,19] = 0.004818030126325816
$data[48,0] = 'ibes_1|fwdepsqcut-46|dense2｜new with indi code -fix space'
$data[48,1] = 0.009301347310769852
$data[48,2] = 0.008733878054388279
$data[48,3] = 0.0001883109214593075
$data[48,4] = 0.0001961492659581503
$data[48,5] = 0.1628933562675233
$data[48,6] = 0.1280492261183164
$data[48,7] = 0.2582640092197721
$data[48,8] = 14156
$data[48,10] = 0.0004026104577395279
$data[48,12] = 0.0003384092866814527
$data[48,14] = 0.01035467890006862
$data[48,15] = 0.005399385735558927
$data[48,16] = 0.009674189396799987
$data[48,18] = 0.1175458874123242
$data[48,19] = 0.004818030126325817
$data[49,0] = 'ibes_qoq_1|fwdepsqcut|q_1｜ibes_qoq_tune10_2'
$data[49,1] = 0.00303210497017975
$data[49,2] = 0.002499571894821071
$data[49,3] = 0.00001909735738640315
$data[49,4] = 0.00001967007435869721
$data[49,5] = -0.0120856426831728
$data[49,6] = -0.04243741404356438
$data[49,7] = 0.2133795918767654
$data[49,8] = 11611
$data[49,10] = 0.0002080153046483946
$data[49,12] = 0.0001631540263502446
$data[49,14] = 0.005090949099942899
$data[49,15] = 0.001550051802718652
$data[49,16] = 0.004580501696160511
$data[49,18] = -0.002911711704464803
$data[49,19] = 0.001592932311570713
$data[50,0] = 'ibes_1|fwdepsqcut-46|dense2｜top15_infwd_mini -code 0 -exclude_fwd False'
$data[50,1] = 0.008142928261746932
$data[50,2] = 0.008733878054388279
$data[50,3] = 0.0001648112327367591
$data[50,4] = 0.0001961492659581503
$data[50,5] = 0.2673575339309593
$data[50,6] = 0.1280492261183164
$data[50,7] = 0.2582640092197721
$data[50,8] = 14156
$data[50,10] = 0.0003438074749231005
$data[50,12] = 0.0003384092866814527
$data[50,14] = 0.009263751492683912
$data[50,15] = 0.004666820467017538
$data[50,16] = 0.009674189396799987
$data[50,18] = 0.2464320924804251
$data[50,19] = 0.004818030126325817
$data[51,0] = 'ibes_1|ni|rnn_top｜top15_lgbm'
$data[51,1] = 0.01025436961368844
$data[51,2] = 0.01023301201875394
$data[51,3] = 0.0002841844015704704
$data[51,4] = 0.0003551633527675986
$data[51,5] = 0.3369590471682379
$data[51,6] = 0.1713554771880856
$data[51,7] = 0.4827027158358849
$data[51,8] = 1267
$data[51,10] = 0.0006028948928221208
$data[51,12] = 0.0004283113735138171
$data[51,14] = 0.0102050584404969
$data[51,15] = 0.005379622896887214
$data[51,16] = 0.01026809155717979
$data[51,18] = 0.2718477491393597
$data[51,19] = 0.004670472236637264
$data[52,0] = 'ibes_1|fwdepsqcut-46|dense2｜small_new_config_1 -code 0 -exclude_fwd True'
$data[52,1] = 0.009592505378186515
$data[52,2] = 0.008764961593413305
$data[52,3] = 0.0001961445305465736
$data[52,4] = 0.0001968262288307614
$data[52,5] = 0.122423362108982
$data[52,6] = 0.1193733535942066
$data[52,7] = 0.2253493365624453
$data[52,8] = 7054
$data[52,10] = 0.0003426357597594948
$data[52,12] = 0.0002911643753660776
$data[52,14] = 0.01043843340429874
$data[52,15] = 0.005648699862893431
$data[52,16] = 0.009545579499406483
$data[52,18] = 0.08840833195541187
$data[52,19] = 0.004876716066832584

$rng = $ws.Range("A1:T53")
$rng.Value = $data

$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = $false
